$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D7").Value = -7.4376
$ws.Range("B9").Value = 6.276499999999995
$ws.Range("D12").Value = -7.288400000000004
$ws.Range("B18").Value = 7.525199999999993
$ws.Range("B20").Value = 9.321799999999994
$ws.Range("D26").Value = -8.233600000000001
$ws.Range("B27").Value = 5.910300000000002
$ws.Range("D27").Value = -8.801999999999998
$ws.Range("D29").Value = -7.354499999999999
$ws.Range("D37").Value = -7.625899999999994
$ws.Range("D38").Value = -8.281599999999999
$ws.Range("D51").Value = -8.046999999999997
$ws.Range("D55").Value = -8.636299999999999
$ws.Range("B69").Value = 5.561599999999995
$ws.Range("D69").Value = -7.171799999999995
$ws.Range("D70").Value = -7.445399999999996
$ws.Range("B76").Value = 5.251000000000001
$ws.Range("B82").Value = 6.212600000000001
$ws.Range("D83").Value = -8.909499999999998
$ws.Range("D102").Value = -7.566999999999995
